# Begin adding thyroid pk data sent from Katie Paul-Friedman (Pilari et al. 2017).
$wb = $excel.ActiveWorkbook

$tc = $wb.Worksheets.Item("TissueComp")

# Insert a new row above row 15 (shifts old row15.. down to row16..)
$tc.Rows.Item(15).Insert()

# Populate the new row 15 with the Thyroid / Human / Pilari et al. 2017 entry
$tc.Range("A15").Value = "Thyroid"
$tc.Range("B15").Value = "Human"
$tc.Range("C15").Value = "Pilari et al. 2017"
$tc.Range("D15").Value = 0.535
$tc.Range("E15").Value = 0.089

# Match row-height formatting used by the adjoining rows (14 and 15)
$tc.Rows.Item(14).RowHeight = 15.75
$tc.Rows.Item(15).RowHeight = 15.75

# The engine's row-insert doesn't retarget sheet-scoped defined names that
# point at TissueComp rows past the insertion point, so fix those up by hand
# (mirrors what Excel itself does automatically on a real row insert).
$wb.Names.Item("TissueComp!bbib18").RefersTo = "=TissueComp!`$A`$32"
$wb.Names.Item("TissueComp!bbib8").RefersTo = "=TissueComp!`$A`$33"
$wb.Names.Item("VolumeFlow!bbib8").RefersTo = "=TissueComp!`$A`$37"
$wb.Names.Item("TissueComp!tblfn10").RefersTo = "=TissueComp!`$A`$31"
$wb.Names.Item("TissueComp!tblfn8").RefersTo = "=TissueComp!`$A`$29"
$wb.Names.Item("TissueComp!tblfn9").RefersTo = "=TissueComp!`$A`$30"

# Update view state: TissueComp becomes the active/visible sheet
$tc.Activate()
$excel.ActiveWindow.ScrollRow = 7
$tc.Range("F15").Select()

# VolumeFlow: scroll view only
$vf = $wb.Worksheets.Item("VolumeFlow")
$vf.Activate()
$excel.ActiveWindow.ScrollRow = 25
$vf.Range("A50:XFD61").Select()

# Basic PK: selection moves, no longer the active tab
$bpk = $wb.Worksheets.Item("Basic PK")
$bpk.Activate()
$bpk.Range("A18").Select()

# Percent BW: selection changes
$pbw = $wb.Worksheets.Item("Percent BW")
$pbw.Activate()
$pbw.Range("O19").Select()

# Leave TissueComp as the final active sheet/tab
$tc.Activate()
$tc.Range("F15").Select()
